$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 11 scores (student "Nguyễn Đắc Huy")
$ws.Range("E11").Value = 7
$ws.Range("H11").Value = 8
$ws.Range("I11").Value = 6.5
$ws.Range("J11").Value = "C+"
